$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Title line: "Alex (Sean) Wall and Noah Stinson" -> "Sean (Alex) Wall and
#    Noah Stinson", materialised as two runs ("Sean (Alex) Wall" /
#    " and Noah Stinson") with identical run formatting, matching the target
#    markup. A temporary bookmark inserted at the split point (then removed)
#    forces the run boundary without altering any formatting.
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs(3).Range
$titlePara.Find.Execute("Alex (Sean) Wall and Noah Stinson", $true, $false, `
    $false, $false, $false, $true, 1, $false, `
    "Sean (Alex) Wall and Noah Stinson", 2) | Out-Null

$titlePara2 = $d.Paragraphs(3).Range
$titleScan = $d.Range($titlePara2.Start, $titlePara2.End)
$titleScan.Find.Execute("Sean (Alex) Wall") | Out-Null
$splitPoint = $d.Range($titleScan.End, $titleScan.End)
$d.Bookmarks.Add("TempRunSplit", $splitPoint) | Out-Null
$d.Bookmarks("TempRunSplit").Delete() | Out-Null

# ---------------------------------------------------------------------------
# 2) Drop the stray "_GoBack" bookmark that used to sit in the blank
#    paragraph right after the title (it gets relocated below in step 4).
# ---------------------------------------------------------------------------
$tbl = $d.Tables(1)

# ---------------------------------------------------------------------------
# 3) Table row for task 2 ("Lead By" column): simple text swap within a
#    single run.
# ---------------------------------------------------------------------------
$cellA = $tbl.Cell(2, 2).Range
$rangeA = $d.Range($cellA.Start, $cellA.End)
$rangeA.Find.Execute("Alex (Sean) Wall", $true, $false, $false, $false, `
    $false, $true, 1, $false, "Sean (Alex) Wall", 1) | Out-Null

# ---------------------------------------------------------------------------
# 4) Table row for task 4 ("Lead By" column): text swap, then split into
#    "Sean (Alex)" / " Wall" with the "_GoBack" bookmark re-homed to the
#    split point (matching where Word's last-edit marker ends up).
# ---------------------------------------------------------------------------
$tbl2 = $d.Tables(1)
$cellB = $tbl2.Cell(4, 2).Range
$rangeB = $d.Range($cellB.Start, $cellB.End)
$rangeB.Find.Execute("Alex (Sean) Wall", $true, $false, $false, $false, `
    $false, $true, 1, $false, "Sean (Alex) Wall", 1) | Out-Null

$tbl3 = $d.Tables(1)
$cellB2 = $tbl3.Cell(4, 2).Range
$rangeB2 = $d.Range($cellB2.Start, $cellB2.End)
$rangeB2.Find.Execute("Sean (Alex)") | Out-Null
$goBackPoint = $d.Range($rangeB2.End, $rangeB2.End)
$d.Bookmarks.Add("_GoBack", $goBackPoint) | Out-Null
